$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 112328247
$ws.Range("B5").Value = 77797
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 1249
$ws.Range("F5").Value = "Norsk näverlav"
$ws.Range("G5").Value = "Platismatia norvegica"
$ws.Range("H5").Value = "(Lynge) W.L.Culb. & C.F.Culb."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = ""
$ws.Range("I5").ClearFormats()
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = ""
$ws.Range("K5").ClearFormats()
$ws.Range("P5").Value = "Framnäs, Jmt"
$ws.Range("Q5").Value = 457353
$ws.Range("R5").Value = 7151591
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Strömsund"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Frostviken"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-26"
$ws.Range("Y5").ClearFormats()
$ws.Range("Z5").Value = "12:59"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-26"
$ws.Range("AA5").ClearFormats()
$ws.Range("AB5").Value = "12:59"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").NumberFormat = "@"
$ws.Range("AT5").Value = ""
$ws.Range("AT5").ClearFormats()
$ws.Range("AW5").Value = "Henrik Tykosson"
$ws.Range("AX5").Value = "Henrik Tykosson"
$ws.Range("AY5").NumberFormat = "@"
$ws.Range("AY5").Value = ""
$ws.Range("AY5").ClearFormats()
